# V1.5a production changes, change crystal to one with higher load abilitys
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- X1 oscillator row (row 48): swap crystal for one with higher load ability ---
$ws.Range("D48").Value = "Oscillator XO91 4Pin 7.0x5.0mm"
$ws.Range("F48").Value = "CB3LV-5I-40M0000"
$ws.Range("H48").Value = ""

# --- Footer: bump PCBWay version + production date ---
$ws.Range("A50").Value = "PCBWay Friendly Version v1.6"
$ws.Range("A55").Value = "Date: 05.02.2025 (DD.MM.YYYY)"
